# Natmi following Dr Hou advice
# Update the Mfge8-Itgb3 LR-pair metrics (rows 2-10, columns E:T) with the
# recalculated values (ligand/receptor expressing-cell counts corrected
# from 1 to 3, with all downstream totals/specificities/weights updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(3, 1, 21.66138833333333, 64.98416499999999, 0.2252765553546639, 0.2252765553546639, 3, 1, 8.970048, 26.910144, 0.487108783009476, 0.4871087830094759, 194.30369309664, 1748.73323786976, 0.1097341887193772, 0.1097341887193772)
    3  = @(3, 1, 21.66138833333333, 64.98416499999999, 0.2252765553546639, 0.2252765553546639, 3, 1, 9.012070666666666, 27.036212, 0.489390778604016, 0.489390778604016, 195.2139623981089, 1756.92566158298, 0.1102482688262497, 0.1102482688262497)
    4  = @(3, 1, 21.66138833333333, 64.98416499999999, 0.2252765553546639, 0.2252765553546639, 3, 1, 0.4327576666666667, 1.298273, 0.02350043838650813, 0.02350043838650813, 9.374131871893887, 84.36718684704499, 0.005294097809037066, 0.005294097809037066)
    5  = @(3, 1, 24.68088566666666, 74.04265699999999, 0.2566790650963491, 0.2566790650963491, 3, 1, 8.970048, 26.910144, 0.487108783009476, 0.4871087830094759, 221.388729112512, 1992.498562012608, 0.1250306270230927, 0.1250306270230927)
    6  = @(3, 1, 24.68088566666666, 74.04265699999999, 0.2566790650963491, 0.2566790650963491, 3, 1, 9.012070666666666, 27.036212, 0.489390778604016, 0.489390778604016, 222.4258857439204, 2001.832971695284, 0.1256163675188532, 0.1256163675188532)
    7  = @(3, 1, 24.68088566666666, 74.04265699999999, 0.2566790650963491, 0.2566790650963491, 3, 1, 0.4327576666666667, 1.298273, 0.02350043838650813, 0.02350043838650813, 10.68084249237344, 96.12758243136099, 0.00603207055440326, 0.006032070554403262)
    8  = @(3, 1, 49.812376, 149.437128, 0.518044379548987, 0.5180443795489871, 3, 1, 8.970048, 26.910144, 0.487108783009476, 0.4871087830094759, 446.819403714048, 4021.374633426432, 0.2523439672670061, 0.2523439672670061)
    9  = @(3, 1, 49.812376, 149.437128, 0.518044379548987, 0.5180443795489871, 3, 1, 9.012070666666666, 27.036212, 0.489390778604016, 0.489390778604016, 448.9126525865706, 4040.213873279136, 0.2535261422589131, 0.2535261422589132)
    10 = @(3, 1, 49.812376, 149.437128, 0.518044379548987, 0.5180443795489871, 3, 1, 0.4327576666666667, 1.298273, 0.02350043838650813, 0.02350043838650813, 21.55668760888267, 194.010188479944, 0.0121742700230678, 0.0121742700230678)
}

$cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value2 = $rowVals[$i]
    }
}
